# =========================================================================
# Applies the "adding submitted channel paper 20220324" edit to the CV.
# =========================================================================
$d = $word.ActiveDocument

function Find-And-Get-Range {
    param($paraIndex, [string]$searchText)
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $ok = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Find failed for paragraph $paraIndex : $searchText"
    }
    return $r
}

# -------------------------------------------------------------------------
# 1) Paragraph 31 : "...Zhonghua Zhao, "Separation...," Geophysics (accepted)."
#    -> "...Zhonghua Zhao, (2022), "Separation...," Geophysics (Accepted)."
# -------------------------------------------------------------------------
$r = Find-And-Get-Range 31 "Zhonghua Zhao, "
$insPos = $r.End
$d.Range($insPos, $insPos).Text = "(2022), "

$r = Find-And-Get-Range 31 "(accepted)"
$aStart = $r.Start + 1
$d.Range($aStart, $aStart + 1).Text = "A"

# -------------------------------------------------------------------------
# 2) Paragraph 32 : "...Wenchao Chen, "Eliminating...," Geophysics (accepted)."
#    -> "...Wenchao Chen, (2022), "Eliminating...," Geophysics (Accepted)."
# -------------------------------------------------------------------------
$r = Find-And-Get-Range 32 "Wenchao Chen, ""Eliminating"
$insPos = $r.Start + ("Wenchao Chen, ").Length
$d.Range($insPos, $insPos).Text = "(2022), "

$r = Find-And-Get-Range 32 "(accepted)"
$aStart = $r.Start + 1
$d.Range($aStart, $aStart + 1).Text = "A"

# -------------------------------------------------------------------------
# 3) Paragraph 33 : "...Interpretation, 10: SA59-SA67." -> ", 10, SA59-SA67."
# -------------------------------------------------------------------------
$r = Find-And-Get-Range 33 "Interpretation, 10: SA59-SA67."
$r.Text = "Interpretation, 10, SA59-SA67."

# -------------------------------------------------------------------------
# 4) Paragraph 35 : "...Geophysics, 86: V509-V523." -> ", 86, V509-V523."
# -------------------------------------------------------------------------
$r = Find-And-Get-Range 35 "Geophysics, 86: V509-V523."
$r.Text = "Geophysics, 86, V509-V523."

# -------------------------------------------------------------------------
# 5) Paragraph 38 : "...pp. 1598-1629, March 2020." -> "...pp. 1598-1629, 2020."
# -------------------------------------------------------------------------
$r = Find-And-Get-Range 38 "pp. 1598-1629, March 2020."
$r.Text = "pp. 1598-1629, 2020."

# -------------------------------------------------------------------------
# 6) Paragraph 41 : "...(Under review)" -> "...(Under major revision)"
#    with "major" bolded.
# -------------------------------------------------------------------------
$r = Find-And-Get-Range 41 "Under review)"
$r.Text = "Under major revision)"

$r = Find-And-Get-Range 41 "major"
$r.Font.Bold = 1

# -------------------------------------------------------------------------
# 7) New paragraph after paragraph 41 (same numbered list, numId 37):
#    "Dawei Liu, Wei Wang, Xiaokai Wang, Zhensheng Shi, Mauricio D. Sacchi
#    ,Wenchao Chen, (2022), "Improving sparse representation with deep
#    learning: a workflow for separating strong background interference,"
#    Geophysics. (With journal)."
# -------------------------------------------------------------------------
$p41 = $d.Paragraphs(41)
$p41.Range.InsertParagraphAfter()

$p42 = $d.Paragraphs(42)
$insPos = $p42.Range.Start
$d.Range($insPos, $insPos).Text = "Dawei Liu"
$boldRange = $d.Range($insPos, $insPos + ("Dawei Liu").Length)
$boldRange.Font.Bold = 1

$p42 = $d.Paragraphs(42)
$insPos = $p42.Range.Start + ("Dawei Liu").Length
$rest = ", Wei Wang, Xiaokai Wang, Zhensheng Shi, Mauricio D. Sacchi ,Wenchao Chen, (2022), ""Improving sparse representation with deep learning: a workflow for separating strong background interference,"" Geophysics. (With journal)."
$d.Range($insPos, $insPos).Text = $rest
$restRange = $d.Range($insPos, $insPos + $rest.Length)
$restRange.Font.Bold = 0

Write-Output "done"
